## Añadido metodo ocr multiple
## Splits the combined upper/lower-case "W/w" weekly regex rows in the
## date_format sheet into separate rows: one that keeps the upper-case-only
## regex (format code unchanged) and a brand new row with the lower-case
## "w" regex + lower-case format code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("date_format")

# --- 1. Temporarily drop the trailing "LargeMonth" rows so the sheet's
#        used range shrinks back to 39 rows -----------------------------
$ws.Range("A40:E41").EntireRow.Delete()

# --- 2. Re-apply the AutoFilter so its stored range becomes exactly
#        A1:E40 (one row further than the original A1:E39) -------------
$ws.AutoFilterMode = $false
$ws.Range("A1:E40").AutoFilter()

# --- 3. Restore the two "LargeMonth" rows (format + values) back as
#        rows 40 and 41 --------------------------------------------------
$ws.Range("A38:E38").Copy()
$ws.Range("A40:E41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(40, 1).Value = "LargeMonth dd, yyyy"
$ws.Cells.Item(40, 2).Value = "December 31, 2012"
$ws.Cells.Item(40, 3).Value = "[a-zA-Z]+ \d{1,2}, *\d{3,4}"
$ws.Cells.Item(40, 4).Value = "%B %d, %Y"
$ws.Cells.Item(41, 1).Value = "LargeMonth dd, yyyy"
$ws.Cells.Item(41, 2).Value = "December 31, 2012"
$ws.Cells.Item(41, 3).Value = "[a-zA-Z]+ \d{1,2}, *\d{2}"
$ws.Cells.Item(41, 4).Value = "%B %d, %y"

# --- 4. Fix up the two "combined" W/w rows so they only match the
#        upper-case "W" variant any more --------------------------------
$ws.Cells.Item(38, 3).Value = "W \d{1,2}\.\d{3,4}"
$ws.Cells.Item(39, 3).Value = "W \d{1,2}\.\d{2}"

# --- 5. Insert a brand new row right after row 38 for the lower-case
#        "w \d{1,2}\.\d{3,4}" variant ------------------------------------
$ws.Rows.Item(39).Insert()
$ws.Range("A38:E38").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)
$ws.Cells.Item(39, 1).Value = "W ww.yyyy"
$ws.Cells.Item(39, 2).Value = "W 52.2012"
$ws.Cells.Item(39, 3).Value = "w \d{1,2}\.\d{3,4}"
$ws.Cells.Item(39, 4).Value = "w %W.%Y"
$ws.Cells.Item(39, 5).Value = 'Añadir "Monday" al text input y "%A" al format code para que coja el primer día de la semana'

# --- 6. Insert a second new row right after row 40 (the other fixed-up
#        combined row, now shifted down to row 40) for the lower-case
#        "w \d{1,2}\.\d{2}" variant --------------------------------------
$ws.Rows.Item(41).Insert()
$ws.Range("A40:E40").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)
$ws.Cells.Item(41, 1).Value = "W ww.yyyy"
$ws.Cells.Item(41, 2).Value = "W 52.2012"
$ws.Cells.Item(41, 3).Value = "w \d{1,2}\.\d{2}"
$ws.Cells.Item(41, 4).Value = "w %W.%y"
$ws.Cells.Item(41, 5).Value = 'Añadir "Monday" al text input y "%A" al format code para que coja el primer día de la semana'

$excel.CutCopyMode = $false

# --- 7. Re-point the sheet's hidden _FilterDatabase defined name --------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "date_format!_FilterDatabase") {
        $n.RefersTo = "=date_format!`$A`$1:`$E`$40"
    }
}

# --- 8. Update the sheet's selection to match the author's last position -
$ws.Range("C41").Select()
